$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column Z header is renamed first so its new shared string ("entuty_status_id")
# is appended to the shared-strings table before column Y's ("added_on_entuity"),
# matching the author's edit order (entuty_status_id=60, added_on_entuity=61).
#   Z1: "Entuty_add-id"   -> "entuty_status_id"
#   Y1: "entuity_status"  -> "added_on_entuity"
$ws.Range("Z1").Value = "entuty_status_id"
$ws.Range("Y1").Value = "added_on_entuity"

# Y2 data cell changes from the numeric 1 to the text "Yes/No"
$ws.Range("Y2").Value = "Yes/No"

# Update the view: scroll back to show column A and move the selection to F2
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F2").Select()
